$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06747606590576538
$ws.Range("D2").Value = 0.1186082610989203
$ws.Range("E2").Value = 0.1294418084900464
$ws.Range("F2").Value = 1.916543709514201
$ws.Range("G2").Value = 1.275509177230489
$ws.Range("H2").Value = 1.201867251022534
$ws.Range("I2").Value = 0.5977971328830698
$ws.Range("J2").Value = 0.1679489714918105
$ws.Range("K2").Value = 0.9336118786398231
$ws.Range("M2").Value = 0.3570140207309791
$ws.Range("N2").Value = 1.827886918869599

$ws.Range("B3").Value = 0.05917170631208535
$ws.Range("D3").Value = 0.1173768301360809
$ws.Range("E3").Value = 0.1290353841054745
$ws.Range("F3").Value = 1.912345820883672
$ws.Range("G3").Value = 1.269815302266551
$ws.Range("H3").Value = 1.204628491038235
$ws.Range("I3").Value = 0.6040713770614285
$ws.Range("J3").Value = 0.1680844295775081
$ws.Range("K3").Value = 0.8545992344685658
$ws.Range("M3").Value = 0.3399944152535355
$ws.Range("N3").Value = 1.847432209408699

$ws.Range("B4").Value = 0.054067404551148
$ws.Range("D4").Value = 0.1166581953051775
$ws.Range("E4").Value = 0.1288370127842953
$ws.Range("F4").Value = 1.910801666421847
$ws.Range("G4").Value = 1.26711460737495
$ws.Range("H4").Value = 1.206888502574671
$ws.Range("I4").Value = 0.6081821138498107
$ws.Range("J4").Value = 0.1682418537829946
$ws.Range("K4").Value = 0.8063986839827351
$ws.Range("M4").Value = 0.3297090373926324
$ws.Range("N4").Value = 1.860035993907022

$ws.Range("B5").Value = 0.05198617059970445
$ws.Range("D5").Value = 0.1163748117610481
$ws.Range("E5").Value = 0.1287690655726479
$ws.Range("F5").Value = 1.91043215119754
$ws.Range("G5").Value = 1.266213774609568
$ws.Range("H5").Value = 1.207951445479182
$ws.Range("I5").Value = 0.6099222312732646
$ws.Range("J5").Value = 0.168324684643455
$ws.Range("K5").Value = 0.7868357652513964
$ws.Range("M5").Value = 0.3255591980950072
$ws.Range("N5").Value = 1.865323715095226

$ws.Range("B6").Value = 0.0516405163377982
$ws.Range("D6").Value = 0.1163283291641051
$ws.Range("E6").Value = 0.1287585622187102
$ws.Range("F6").Value = 1.910386478079857
$ws.Range("G6").Value = 1.266076246646207
$ws.Range("H6").Value = 1.208136520926203
$ws.Range("I6").Value = 0.6102150998256288
$ws.Range("J6").Value = 0.1683395670597641
$ws.Range("K6").Value = 0.7835921514787003
$ws.Range("M6").Value = 0.3248726324719513
$ws.Range("N6").Value = 1.866210891297825

$ws.Range("B7").Value = 0.05403934089844142
$ws.Range("D7").Value = 0.1166543351093381
$ws.Range("E7").Value = 0.1288360442004652
$ws.Range("F7").Value = 1.910795631479203
$ws.Range("G7").Value = 1.26710165013418
$ws.Range("H7").Value = 1.206902262971951
$ws.Range("I7").Value = 0.6082053186659433
$ws.Range("J7").Value = 0.168242895226264
$ws.Range("K7").Value = 0.8061345304788574
$ws.Range("M7").Value = 0.3296529028812643
$ws.Range("N7").Value = 1.860106692274441

$ws.Range("B8").Value = 0.0646139387360023
$ws.Range("D8").Value = 0.11817591291382
$ws.Range("E8").Value = 0.1292910649967247
$ws.Range("F8").Value = 1.914881729974837
$ws.Range("G8").Value = 1.273380676556144
$ws.Range("H8").Value = 1.202702139460385
$ws.Range("I8").Value = 0.5999068777149468
$ws.Range("J8").Value = 0.1679802706573135
$ws.Range("K8").Value = 0.9063034324389889
$ws.Range("M8").Value = 0.3511115329633014
$ws.Range("N8").Value = 1.83450098943326

$ws.Range("B9").Value = 0.08530160974221701
$ws.Range("D9").Value = 0.1214553235062894
$ws.Range("E9").Value = 0.1305886040058581
$ws.Range("F9").Value = 1.931101884618741
$ws.Range("G9").Value = 1.29201988525061
$ws.Range("H9").Value = 1.198946845800592
$ws.Range("I9").Value = 0.5856831965153386
$ws.Range("J9").Value = 0.1680542541174006
$ws.Range("K9").Value = 1.105217175812555
$ws.Range("M9").Value = 0.3944964668745996
$ws.Range("N9").Value = 1.789072416290784

$ws.Range("B10").Value = 0.1004639098737954
$ws.Range("D10").Value = 0.1240429912801417
$ws.Range("E10").Value = 0.1317881343001801
$ws.Range("F10").Value = 1.948037920843106
$ws.Range("G10").Value = 1.309595371101068
$ws.Range("H10").Value = 1.198923117874187
$ws.Range("I10").Value = 0.576482108764246
$ws.Range("J10").Value = 0.168467705032505
$ws.Range("K10").Value = 1.252882035373204
$ws.Range("M10").Value = 0.4271672388650884
$ws.Range("N10").Value = 1.758612019585579

$ws.Range("B11").Value = 0.1073522184633759
$ws.Range("D11").Value = 0.1252585192878897
$ws.Range("E11").Value = 0.1323871438504121
$ws.Range("F11").Value = 1.956836137068308
$ws.Range("G11").Value = 1.318439290586326
$ws.Range("H11").Value = 1.1995071633294
$ws.Range("I11").Value = 0.5725674342309439
$ws.Range("J11").Value = 0.1687337975986054
$ws.Range("K11").Value = 1.320392330018763
$ws.Range("M11").Value = 0.4422031990955304
$ws.Range("N11").Value = 1.745387869234511

$ws.Range("B12").Value = 0.1099591838242873
$ws.Range("D12").Value = 0.1257242888607806
$ws.Range("E12").Value = 0.1326216253419439
$ws.Range("F12").Value = 1.960325322447375
$ws.Range("G12").Value = 1.321910687886856
$ws.Range("H12").Value = 1.199813914199552
$ws.Range("I12").Value = 0.5711240139991496
$ws.Range("J12").Value = 0.1688457761891939
$ws.Range("K12").Value = 1.346005050584608
$ws.Range("M12").Value = 0.4479218747602829
$ws.Range("N12").Value = 1.74047121107481

$ws.Range("B13").Value = 0.1093977962023587
$ws.Range("D13").Value = 0.1256237342268633
$ws.Range("E13").Value = 0.1325707856613434
$ws.Range("F13").Value = 1.959566856454416
$ws.Range("G13").Value = 1.32115761115142
$ws.Range("H13").Value = 1.199744042678645
$ws.Range("I13").Value = 0.5714331463081379
$ws.Range("J13").Value = 0.168821160844665
$ws.Range("K13").Value = 1.340486762733178
$ws.Range("M13").Value = 0.4466891507771109
$ws.Range("N13").Value = 1.741526050742915

$ws.Range("B14").Value = 0.1075667259130739
$ws.Range("D14").Value = 0.1252967289462958
$ws.Range("E14").Value = 0.1324062815903666
$ws.Range("F14").Value = 1.957120037000905
$ws.Range("G14").Value = 1.318722429509108
$ws.Range("H14").Value = 1.199530684362713
$ws.Range("I14").Value = 0.5724479018039794
$ws.Range("J14").Value = 0.1687427853930359
$ws.Range("K14").Value = 1.322498543470886
$ws.Range("M14").Value = 0.4426731799432488
$ws.Range("N14").Value = 1.744981547652105

$ws.Range("B15").Value = 0.1064449434431509
$ws.Range("D15").Value = 0.1250971405324535
$ws.Range("E15").Value = 0.1323065136989818
$ws.Range("F15").Value = 1.955641805101109
$ws.Range("G15").Value = 1.317246760521073
$ws.Range("H15").Value = 1.199411143464829
$ws.Range("I15").Value = 0.5730745460674598
$ws.Range("J15").Value = 0.1686962386221253
$ws.Range("K15").Value = 1.311486483951569
$ws.Range("M15").Value = 0.4402165189954985
$ws.Range("N15").Value = 1.747109998391982

$ws.Range("B16").Value = 0.1000135442948391
$ws.Range("D16").Value = 0.1239643222851754
$ws.Range("E16").Value = 0.1317500591638741
$ws.Range("F16").Value = 1.947484965210705
$ws.Range("G16").Value = 1.309034508797907
$ws.Range("H16").Value = 1.198896922906812
$ws.Range("I16").Value = 0.5767433958130148
$ws.Range("J16").Value = 0.1684518848013852
$ws.Range("K16").Value = 1.248476822408406
$ws.Range("M16").Value = 0.4261880920026471
$ws.Range("N16").Value = 1.759488983360961

$ws.Range("B17").Value = 0.09606563131599444
$ws.Range("D17").Value = 0.1232791747608033
$ws.Range("E17").Value = 0.1314223388125129
$ws.Range("F17").Value = 1.942761321029494
$ws.Range("G17").Value = 1.304214189176605
$ws.Range("H17").Value = 1.198733851661729
$ws.Range("I17").Value = 0.5790635265105344
$ws.Range("J17").Value = 0.1683219601974244
$ws.Range("K17").Value = 1.209908418208215
$ws.Range("M17").Value = 0.4176265529093541
$ws.Range("N17").Value = 1.767245139771091

$ws.Range("B18").Value = 0.09379405055588563
$ws.Range("D18").Value = 0.1228887116316315
$ws.Range("E18").Value = 0.1312388642865798
$ws.Range("F18").Value = 1.940147360205131
$ws.Range("G18").Value = 1.301521535508186
$ws.Range("H18").Value = 1.198696043458526
$ws.Range("I18").Value = 0.5804235058458822
$ws.Range("J18").Value = 0.1682545739881647
$ws.Range("K18").Value = 1.187756617716673
$ws.Range("M18").Value = 0.412718553750608
$ws.Range("N18").Value = 1.771765797707332

$ws.Range("B19").Value = 0.09302479257586072
$ws.Range("D19").Value = 0.1227571298589822
$ws.Range("E19").Value = 0.1311776060401648
$ws.Range("F19").Value = 1.939279994878959
$ws.Range("G19").Value = 1.300623555839451
$ws.Range("H19").Value = 1.19869285647826
$ws.Range("I19").Value = 0.5808883507812581
$ws.Range("J19").Value = 0.168233019431959
$ws.Range("K19").Value = 1.18026185969461
$ws.Range("M19").Value = 0.411059605010216
$ws.Range("N19").Value = 1.773306634867229

$ws.Range("B20").Value = 0.09648598221343718
$ws.Range("D20").Value = 0.123351736005084
$ws.Range("E20").Value = 0.131456705638989
$ws.Range("F20").Value = 1.943253504279653
$ws.Range("G20").Value = 1.304719052004089
$ws.Range("H20").Value = 1.198745416177672
$ws.Range("I20").Value = 0.5788139051780306
$ws.Range("J20").Value = 0.1683350309435454
$ws.Range("K20").Value = 1.214010810084403
$ws.Range("M20").Value = 0.4185362498070404
$ws.Range("N20").Value = 1.766413322689648

$ws.Range("B21").Value = 0.108104597443031
$ws.Range("D21").Value = 0.1253926300718575
$ws.Range("E21").Value = 0.132454392977035
$ws.Range("F21").Value = 1.95783445140286
$ws.Range("G21").Value = 1.319434376778389
$ws.Range("H21").Value = 1.199591029709779
$ws.Range("I21").Value = 0.5721487853916756
$ws.Range("J21").Value = 0.1687655018097658
$ws.Range("K21").Value = 1.327780818928829
$ws.Range("M21").Value = 0.4438520935153321
$ws.Range("N21").Value = 1.743964112841903

$ws.Range("B22").Value = 0.1156892803033713
$ws.Range("D22").Value = 0.1267583673530766
$ws.Range("E22").Value = 0.1331510176120396
$ws.Range("F22").Value = 1.968282014848469
$ws.Range("G22").Value = 1.32976527866839
$ws.Range("H22").Value = 1.200642568900804
$ws.Range("I22").Value = 0.5680199712098428
$ws.Range("J22").Value = 0.1691122084716241
$ws.Range("K22").Value = 1.402415912924766
$ws.Range("M22").Value = 0.4605424492215207
$ws.Range("N22").Value = 1.729822958138532

$ws.Range("B23").Value = 0.1116420524452764
$ws.Range("D23").Value = 0.1260265438150441
$ws.Range("E23").Value = 0.1327751436697966
$ws.Range("F23").Value = 1.962621887669357
$ws.Range("G23").Value = 1.324186072861352
$ws.Range("H23").Value = 1.200035678425479
$ws.Range("I23").Value = 0.5702027989793166
$ws.Range("J23").Value = 0.168921184221297
$ws.Range("K23").Value = 1.362556285133962
$ws.Range("M23").Value = 0.4516212643249347
$ws.Range("N23").Value = 1.737321766354847

$ws.Range("B24").Value = 0.09629594739378433
$ws.Range("D24").Value = 0.1233189203573062
$ws.Range("E24").Value = 0.1314411530195976
$ws.Range("F24").Value = 1.943030671369982
$ws.Range("G24").Value = 1.304490558667794
$ws.Range("H24").Value = 1.198740013593579
$ws.Range("I24").Value = 0.5789266776661002
$ws.Range("J24").Value = 0.1683290988914479
$ws.Range("K24").Value = 1.212156051130876
$ws.Range("M24").Value = 0.4181249317749831
$ws.Range("N24").Value = 1.766789195495764

$ws.Range("B25").Value = 0.0797110276342039
$ws.Range("D25").Value = 0.1205367094084409
$ws.Range("E25").Value = 0.1301942863941257
$ws.Range("F25").Value = 1.92583364308436
$ws.Range("G25").Value = 1.28629769084759
$ws.Range("H25").Value = 1.199482665399159
$ws.Range("I25").Value = 0.5893117051482584
$ws.Range("J25").Value = 0.1679711818124687
$ws.Range("K25").Value = 1.051138983696006
$ws.Range("M25").Value = 0.3826199730602227
$ws.Range("N25").Value = 1.800849946637276
